{"js": "// Replace \"on Windows Vista or Windows 7\" with \"on Windows\" in the\n// \"Participants must also do most of their web browsing using Firefox on\n// Windows Vista or Windows 7.\" sentence, per the commit\n// \"Added Windows support to ethics docs.\"\nconst searchResults = context.document.body.search(\"on Windows Vista or Windows 7\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Target phrase \"on Windows Vista or Windows 7\" not found in document.');\n}\n\nsearchResults.items[0].insertText(\"on Windows\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace \"on Windows Vista or Windows 7\" with \"on Windows\" in the\n# \"Participants must also do most of their web browsing using Firefox on\n# Windows Vista or Windows 7.\" sentence, per the commit\n# \"Added Windows support to ethics docs.\"\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"on Windows Vista or Windows 7\"\n$find.Replacement.Text = \"on Windows\"\n$find.Forward = $true\n$find.Wrap = $wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, [ref]$find.Replacement.Text, $wdReplaceOne) | Out-Null\n"}
